# "#5: property boat&car done"
# The 汽車 (Car) sheet (sheet3, 3rd tab) previously only had six
# data columns (B..G) and its "header" row was actually a stray
# duplicate of the data row. This adds proper column headers and
# extends the row with the remaining common property columns
# (property_category, category, date, legislator_name,
# legislator_id, source_file, index), matching the layout already
# used on the other property sheets.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# --- extend the header row (row 1) -----------------------------------
# Give the new header cells (H1:N1) the same look as the existing
# header cells before filling in their text.
$ws.Range("B1").Copy($ws.Range("H1:N1"))

$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "capacity"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "register_date"
$ws.Range("F1").Value = "register_reason"
$ws.Range("G1").Value = "acquire_value"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# --- extend the data row (row 2) --------------------------------------
# B2:G2 already hold the correct data; only H2:N2 are new.
$ws.Range("H2").Value = "land"
$ws.Range("I2").Value = "normal"

# "2012-03-22" looks like a date, and Excel would otherwise silently
# convert it to a date serial. Force text entry, then drop back to the
# default/Normal style so no stray number format sticks around.
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "2012-03-22"
$ws.Range("J2").Style = "Normal"

$ws.Range("K2").Value = "許添財"
$ws.Range("L2").Value = 639
$ws.Range("M2").Value = "tmpb8d31"
$ws.Range("N2").Value = 32
